$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.047002792358398
$ws.Range("B1").Value = 3.163487434387207
$ws.Range("C1").Value = 5.855867862701416
$ws.Range("D1").Value = 2.347394704818726
$ws.Range("E1").Value = 0.9339359402656555
